$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date header row (row 2) from column AQ (43) through CJ (88)
# with consecutive dates 44044..44089, reusing the same date style as AP2 (s=4)
# by copying its formatting before writing the values.
$ws.Range("AP2").Copy()
$ws.Range("AQ2:CJ2").PasteSpecial(-4122)

$col = 43
$val = 44044
while ($col -le 88) {
    $ws.Cells.Item(2, $col).Value = $val
    $col = $col + 1
    $val = $val + 1
}

# Add a new "3.5 h" entry next to AP26, re-using AP26's style (s=1)
$ws.Range("AP26").Copy()
$ws.Range("AQ26").PasteSpecial(-4122)
$ws.Range("AQ26").Value = "3.5 h"

# Update the active selection shown in the sheet view
$ws.Range("CL31").Select()
